$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 221 (Distributed Energy / Electric boiler / FI00_others_dheat, year 2030) ---
$ws.Cells.Item(221, 1).Value = "FI00"
$ws.Cells.Item(221, 2).Value = "FI00_others_dheat"
$ws.Cells.Item(221, 3).Value = "Electric boiler"
$ws.Cells.Item(221, 4).Value = "Distributed Energy"
$ws.Cells.Item(221, 5).Value = 2030
$ws.Cells.Item(221, 7).Value = 250

# --- Update an existing capacity value (Heat_capa for Electric boiler / 2040 entry) ---
$ws.Cells.Item(142, 7).Value = 350

# --- Re-apply the sheet's AutoFilter: drop the Generator_ID ("Kraft process recovery
#     boiler") criterion, and point the Heatnode criterion at FI00_others_dheat. This
#     also updates which rows are hidden/visible to match the new filter. ---
[void]$ws.Range("A1:J220").AutoFilter(2, @("FI00_others_dheat"))
[void]$ws.Range("A1:J220").AutoFilter(3)

# --- Update the remembered selection ---
[void]$ws.Range("D226").Select()
